$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the contact-info runs (removes the spell-check-split proofErr pair
#    around "danielmartincraig" and appends the linkedin part into one run)
# ---------------------------------------------------------------------------
$contactText = "(803)389-6750 " + [char]0x2022 + " danielmartincraig@gmail.com " + [char]0x2022 + " github.com/danielmartincraig " + [char]0x2022 + " linkedin.com/danielcraig23"
$rngContact = $d.Content.Duplicate
$rngContact.Find.Execute("github.com/danielmartincraig " + [char]0x2022 + " linkedin.com/danielcraig23", $true, $false, $false, $false, $false, $true, 1, $false, "github.com/danielmartincraig " + [char]0x2022 + " linkedin.com/danielcraig23", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Merge the "Web Engineering I and II" runs (removes the gramStart/gramEnd
#    proofErr pair)
# ---------------------------------------------------------------------------
$rngWeb = $d.Content.Duplicate
$webText = [char]0x2022 + "    Web Engineering I and II"
$rngWeb.Find.Execute($webText, $true, $false, $false, $false, $false, $true, 1, $false, $webText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new "OBJECTIVE:" paragraph right after the contact-info line
#    (paragraph 3), before "EDUCATION:". Built to match the target run
#    structure: a plain "OBJECTIVE: " run (inherits Heading1 style formatting)
#    followed by three sz=24 runs, with the _GoBack bookmark sitting between
#    "TRC" and " on a full-time basis".
# ---------------------------------------------------------------------------
$pContact = $d.Paragraphs.Item(3)
$pContact.Range.InsertParagraphAfter()
$pObjective = $d.Paragraphs.Item(4)

$objectiveXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
  + '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r>' `
  + '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Eager to drive back-end solutions at </w:t></w:r>' `
  + '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>TRC</w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
  + '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> on a full-time basis</w:t></w:r>' `
  + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pObjective.Range.InsertXML($objectiveXml)

# ---------------------------------------------------------------------------
# 4. Remove the stray _GoBack bookmark that used to sit at the very end of
#    the document (after "Fluent in Spanish") - it has now moved up into the
#    new OBJECTIVE paragraph above, so the old one must go. We rebuild the
#    last paragraph (preserving its three existing runs verbatim) just
#    before the old bookmark-carrying paragraph, then delete the old one.
# ---------------------------------------------------------------------------
$totalParas = $d.Paragraphs.Count
$pBeforeLast = $d.Paragraphs.Item($totalParas - 1)
$pBeforeLast.Range.InsertParagraphAfter()
$pNewLast = $d.Paragraphs.Item($totalParas)

$fluentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
  + '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">' + [char]0x2022 + ' </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">  </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> Fluent in Spanish</w:t></w:r>' `
  + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pNewLast.Range.InsertXML($fluentXml)

$oldLastIndex = $d.Paragraphs.Count
$pOldLast = $d.Paragraphs.Item($oldLastIndex)
$pOldLast.Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
